$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4169439
$ws.Range("I40").Value = 3196.6667
$ws.Range("K40").Value = 3196.6667
$ws.Range("M40").Value = -3021.6667
$ws.Range("H53").Value = 777.5
$ws.Range("I53").Value = 1037.4
$ws.Range("K53").Value = 1037.4
$ws.Range("M53").Value = -400.4000000000001
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("H74").Value = 65229820
$ws.Range("I74").Value = 166672960
$ws.Range("K74").Value = 166672960
$ws.Range("M74").Value = -166672024
$ws.Range("H77").Value = 65229820
$ws.Range("I77").Value = 166672960
$ws.Range("K77").Value = 833364800
$ws.Range("M77").Value = -833360120
$ws.Range("H86").Value = 84327336
$ws.Range("I86").Value = 281250660
$ws.Range("J86").Value = 5558007.5
$ws.Range("K86").Value = 281250660
$ws.Range("L86").Value = 5558007.5
$ws.Range("M86").Value = -281249537
$ws.Range("N86").Value = -5560253.5
$ws.Range("H89").Value = 84327336
$ws.Range("I89").Value = 281250660
$ws.Range("J89").Value = 5558007.5
$ws.Range("K89").Value = 1406253300
$ws.Range("L89").Value = 27790037.5
$ws.Range("M89").Value = -1406247684
$ws.Range("N89").Value = -27801269.5
$ws.Range("H92").Value = 27778322
$ws.Range("I92").Value = 466.7931
$ws.Range("J92").Value = 142858000
$ws.Range("K92").Value = 466.7931
$ws.Range("L92").Value = 142858000
$ws.Range("M92").Value = 781.2069
$ws.Range("N92").Value = -142860496
$ws.Range("H103").Value = 1336.1818
$ws.Range("J103").Value = 1357.6666
$ws.Range("L103").Value = 4072.9998
$ws.Range("N103").Value = -5244.9998
$ws.Range("H106").Value = 5381.3335
$ws.Range("I106").Value = 5381.3335
$ws.Range("K106").Value = 5381.3335
$ws.Range("M106").Value = -4750.3335
$ws.Range("H111").Value = 8930386
$ws.Range("J111").Value = 3072.1428
$ws.Range("L111").Value = 9216.428400000001
$ws.Range("N111").Value = -15350.4284
$ws.Range("H127").Value = 1089
$ws.Range("I127").Value = 1089
$ws.Range("K127").Value = 3267
$ws.Range("M127").Value = 1693
$ws.Range("H129").Value = 1412.2413
$ws.Range("I129").Value = 856.5333000000001
$ws.Range("J129").Value = 2007.6428
$ws.Range("K129").Value = 2569.5999
$ws.Range("L129").Value = 6022.928400000001
$ws.Range("M129").Value = 2430.4001
$ws.Range("N129").Value = -16022.9284
$ws.Range("H132").Value = 885.383
$ws.Range("I132").Value = 861.15216
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2583.45648
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -53.45647999999983
$ws.Range("N132").Value = -11060
$ws.Range("H137").Value = 5477.4688
$ws.Range("I137").Value = 3139
$ws.Range("J137").Value = 7540.8237
$ws.Range("K137").Value = 9417
$ws.Range("L137").Value = 22622.4711
$ws.Range("M137").Value = -6867
$ws.Range("N137").Value = -27722.4711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27028206
$ws.Range("I2").Value = 912.3077
$ws.Range("K2").Value = 912.3077
$ws.Range("M2").Value = -799.3077
$ws.Range("H32").Value = 3181232.8
$ws.Range("I32").Value = 3283896.2
$ws.Range("K32").Value = 3283896.2
$ws.Range("M32").Value = -3283609.2
$ws.Range("H45").Value = 5071.6313
$ws.Range("I45").Value = 1756.7142
$ws.Range("J45").Value = 7005.3335
$ws.Range("K45").Value = 1756.7142
$ws.Range("L45").Value = 7005.3335
$ws.Range("M45").Value = -1379.7142
$ws.Range("N45").Value = -7759.3335
$ws.Range("H63").Value = 1584.4286
$ws.Range("I63").Value = 1181.8334
$ws.Range("K63").Value = 1181.8334
$ws.Range("M63").Value = -495.8334
$ws.Range("H66").Value = 1584.4286
$ws.Range("I66").Value = 1181.8334
$ws.Range("K66").Value = 5909.166999999999
$ws.Range("M66").Value = -2477.166999999999
$ws.Range("H116").Value = 27028206
$ws.Range("I116").Value = 912.3077
$ws.Range("K116").Value = 912.3077
$ws.Range("M116").Value = 1381.6923
$ws.Range("H132").Value = 5946.091
$ws.Range("I132").Value = 1572.1177
$ws.Range("K132").Value = 4716.3531
$ws.Range("M132").Value = -2186.3531

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27028206
$ws.Range("I3").Value = 912.3077
$ws.Range("K3").Value = 912.3077
$ws.Range("M3").Value = -798.3077
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = ""
$ws.Range("H64").Value = 11906238
$ws.Range("I64").Value = 30303876
$ws.Range("J64").Value = 1884.2354
$ws.Range("K64").Value = 30303876
$ws.Range("L64").Value = 1884.2354
$ws.Range("M64").Value = -30303651
$ws.Range("N64").Value = -2334.2354
$ws.Range("H67").Value = 11906238
$ws.Range("I67").Value = 30303876
$ws.Range("J67").Value = 1884.2354
$ws.Range("K67").Value = 30303876
$ws.Range("L67").Value = 1884.2354
$ws.Range("M67").Value = -30303096
$ws.Range("N67").Value = -3444.2354
$ws.Range("H80").Value = 20833720
$ws.Range("J80").Value = 412.3846
$ws.Range("L80").Value = 412.3846
$ws.Range("N80").Value = -2408.3846
$ws.Range("H82").Value = 14128
$ws.Range("I82").Value = 14128
$ws.Range("K82").Value = 14128
$ws.Range("M82").Value = -13745
$ws.Range("H83").Value = 20833720
$ws.Range("J83").Value = 412.3846
$ws.Range("L83").Value = 2061.923
$ws.Range("N83").Value = -12045.923
$ws.Range("H85").Value = 14128
$ws.Range("I85").Value = 14128
$ws.Range("K85").Value = 14128
$ws.Range("M85").Value = -12802
$ws.Range("H105").Value = 4496.273
$ws.Range("I105").Value = 3362
$ws.Range("K105").Value = 3362
$ws.Range("M105").Value = -1615
$ws.Range("H107").Value = 93751230
$ws.Range("I107").Value = 102273896
$ws.Range("K107").Value = 102273896
$ws.Range("M107").Value = -102271976
$ws.Range("H124").Value = 50026
$ws.Range("J124").Value = 50026
$ws.Range("L124").Value = 50026
$ws.Range("N124").Value = -59846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 901
$ws.Range("I6").Value = 901
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 901
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -788
$ws.Range("N6").Value = ""
$ws.Range("H16").Value = 6653.231
$ws.Range("I16").Value = 3747.5
$ws.Range("J16").Value = 7944.6665
$ws.Range("K16").Value = 3747.5
$ws.Range("L16").Value = 7944.6665
$ws.Range("M16").Value = -3460.5
$ws.Range("N16").Value = -8518.666499999999
$ws.Range("H31").Value = 6925.784
$ws.Range("I31").Value = 2273.0833
$ws.Range("J31").Value = 11061.519
$ws.Range("K31").Value = 2273.0833
$ws.Range("L31").Value = 11061.519
$ws.Range("M31").Value = -1978.0833
$ws.Range("N31").Value = -11651.519
$ws.Range("H34").Value = 6925.784
$ws.Range("I34").Value = 2273.0833
$ws.Range("J34").Value = 11061.519
$ws.Range("K34").Value = 2273.0833
$ws.Range("L34").Value = 11061.519
$ws.Range("M34").Value = -2071.0833
$ws.Range("N34").Value = -11465.519
$ws.Range("H68").Value = 59088.25
$ws.Range("J68").Value = 59088.25
$ws.Range("L68").Value = 59088.25
$ws.Range("N68").Value = -60586.25
$ws.Range("H71").Value = 59088.25
$ws.Range("J71").Value = 59088.25
$ws.Range("L71").Value = 177264.75
$ws.Range("N71").Value = -184752.75
$ws.Range("H74").Value = 360330.25
$ws.Range("J74").Value = 472107
$ws.Range("L74").Value = 472107
$ws.Range("N74").Value = -473855
$ws.Range("H77").Value = 360330.25
$ws.Range("J77").Value = 472107
$ws.Range("L77").Value = 1416321
$ws.Range("N77").Value = -1425057
$ws.Range("H99").Value = 3582.3572
$ws.Range("I99").Value = 2332.923
$ws.Range("J99").Value = 4665.2
$ws.Range("K99").Value = 2332.923
$ws.Range("L99").Value = 4665.2
$ws.Range("M99").Value = -834.9229999999998
$ws.Range("N99").Value = -7661.2
$ws.Range("H105").Value = 5495865
$ws.Range("I105").Value = 7143674.5
$ws.Range("J105").Value = 3165.3333
$ws.Range("K105").Value = 7143674.5
$ws.Range("L105").Value = 3165.3333
$ws.Range("M105").Value = -7141927.5
$ws.Range("N105").Value = -6659.3333
$ws.Range("H113").Value = 6653.231
$ws.Range("I113").Value = 3747.5
$ws.Range("J113").Value = 7944.6665
$ws.Range("K113").Value = 3747.5
$ws.Range("L113").Value = 7944.6665
$ws.Range("M113").Value = -1577.5
$ws.Range("N113").Value = -12284.6665
$ws.Range("H126").Value = 3582.3572
$ws.Range("I126").Value = 2332.923
$ws.Range("J126").Value = 4665.2
$ws.Range("K126").Value = 6998.768999999999
$ws.Range("L126").Value = 13995.6
$ws.Range("M126").Value = -4528.768999999999
$ws.Range("N126").Value = -18935.6
$ws.Range("H132").Value = 6293.7
$ws.Range("I132").Value = 3968.35
$ws.Range("J132").Value = 8619.049999999999
$ws.Range("K132").Value = 11905.05
$ws.Range("L132").Value = 25857.15
$ws.Range("M132").Value = -9375.049999999999
$ws.Range("N132").Value = -30917.15
$ws.Range("H134").Value = 4897.375
$ws.Range("J134").Value = 11156.77
$ws.Range("L134").Value = 33470.31
$ws.Range("N134").Value = -38540.31
$ws.Range("H141").Value = 400022.88
$ws.Range("J141").Value = 400022.88
$ws.Range("L141").Value = 400022.88
$ws.Range("N141").Value = -410382.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4378.963
$ws.Range("I34").Value = 1041
$ws.Range("J34").Value = 4646
$ws.Range("K34").Value = 3123
$ws.Range("L34").Value = 13938
$ws.Range("M34").Value = -3039
$ws.Range("N34").Value = -14106
$ws.Range("H39").Value = 12208.889
$ws.Range("J39").Value = 12554.286
$ws.Range("L39").Value = 37662.858
$ws.Range("N39").Value = -38250.858
$ws.Range("H68").Value = 23530790
$ws.Range("I68").Value = 66667588
$ws.Range("J68").Value = 14287192
$ws.Range("K68").Value = 200002764
$ws.Range("L68").Value = 42861576
$ws.Range("M68").Value = -200001953
$ws.Range("N68").Value = -42863198
$ws.Range("H71").Value = 23530790
$ws.Range("I71").Value = 66667588
$ws.Range("J71").Value = 14287192
$ws.Range("K71").Value = 600008292
$ws.Range("L71").Value = 128584728
$ws.Range("M71").Value = -600004236
$ws.Range("N71").Value = -128592840
$ws.Range("H87").Value = 789.3333
$ws.Range("I87").Value = 789.3333
$ws.Range("K87").Value = 2367.9999
$ws.Range("M87").Value = -1119.9999
$ws.Range("H90").Value = 789.3333
$ws.Range("I90").Value = 789.3333
$ws.Range("K90").Value = 7103.9997
$ws.Range("M90").Value = -863.9997000000003
$ws.Range("H97").Value = 404.25
$ws.Range("I97").Value = 300
$ws.Range("K97").Value = 900
$ws.Range("M97").Value = -404
$ws.Range("H116").Value = 2249.5
$ws.Range("I116").Value = 2249.5
$ws.Range("K116").Value = 6748.5
$ws.Range("M116").Value = -3306.5
$ws.Range("H131").Value = 31913.092
$ws.Range("I131").Value = 2000
$ws.Range("J131").Value = 36039.035
$ws.Range("K131").Value = 6000
$ws.Range("L131").Value = 108117.105
$ws.Range("M131").Value = -960
$ws.Range("N131").Value = -118197.105
$ws.Range("H137").Value = 90301.83
$ws.Range("I137").Value = 74253.36
$ws.Range("K137").Value = 222760.08
$ws.Range("M137").Value = -217660.08

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1538724.8
$ws.Range("I2").Value = 129.5
$ws.Range("J2").Value = 2857520.5
$ws.Range("K2").Value = 129.5
$ws.Range("L2").Value = 2857520.5
$ws.Range("M2").Value = -16.5
$ws.Range("N2").Value = -2857746.5
$ws.Range("H80").Value = 2974.75
$ws.Range("I80").Value = 2999.75
$ws.Range("J80").Value = 2949.75
$ws.Range("K80").Value = 2999.75
$ws.Range("L80").Value = 2949.75
$ws.Range("M80").Value = -2001.75
$ws.Range("N80").Value = -4945.75
$ws.Range("H83").Value = 2974.75
$ws.Range("I83").Value = 2999.75
$ws.Range("J83").Value = 2949.75
$ws.Range("K83").Value = 14998.75
$ws.Range("L83").Value = 14748.75
$ws.Range("M83").Value = -10006.75
$ws.Range("N83").Value = -24732.75
$ws.Range("H97").Value = 1060.625
$ws.Range("I97").Value = 628.3333
$ws.Range("J97").Value = 1781.1111
$ws.Range("K97").Value = 628.3333
$ws.Range("L97").Value = 1781.1111
$ws.Range("M97").Value = -132.3333
$ws.Range("N97").Value = -2773.1111
$ws.Range("H102").Value = 2096.8096
$ws.Range("I102").Value = 2145.6667
$ws.Range("K102").Value = 2145.6667
$ws.Range("M102").Value = -523.6667000000002
$ws.Range("H107").Value = 686.2105
$ws.Range("I107").Value = 480
$ws.Range("J107").Value = 724.875
$ws.Range("K107").Value = 480
$ws.Range("L107").Value = 724.875
$ws.Range("M107").Value = 1440
$ws.Range("N107").Value = -4564.875
$ws.Range("H122").Value = 6605404.5
$ws.Range("I122").Value = 6605404.5
$ws.Range("K122").Value = 19816213.5
$ws.Range("M122").Value = -19813763.5
$ws.Range("H126").Value = 5134.643
$ws.Range("I126").Value = 2526.3572
$ws.Range("J126").Value = 7742.9287
$ws.Range("K126").Value = 7579.071599999999
$ws.Range("L126").Value = 23228.7861
$ws.Range("M126").Value = -5109.071599999999
$ws.Range("N126").Value = -28168.7861
$ws.Range("H128").Value = 35249
$ws.Range("J128").Value = 35249
$ws.Range("L128").Value = 35249
$ws.Range("N128").Value = -45209
$ws.Range("H132").Value = 3273.1853
$ws.Range("I132").Value = 1690.0714
$ws.Range("K132").Value = 5070.2142
$ws.Range("M132").Value = -2540.2142
$ws.Range("H134").Value = 45285.57
$ws.Range("J134").Value = 45285.57
$ws.Range("L134").Value = 135856.71
$ws.Range("N134").Value = -140926.71

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 76667
$ws.Range("I2").Value = 65000.5
$ws.Range("J2").Value = 100000
$ws.Range("K2").Value = 65000.5
$ws.Range("L2").Value = 100000
$ws.Range("M2").Value = -64888.5
$ws.Range("N2").Value = -100224
$ws.Range("H7").Value = 3156.0667
$ws.Range("I7").Value = 2535.1365
$ws.Range("J7").Value = 4863.625
$ws.Range("K7").Value = 2535.1365
$ws.Range("L7").Value = 4863.625
$ws.Range("M7").Value = -2423.1365
$ws.Range("N7").Value = -5087.625
$ws.Range("H22").Value = 3326.8
$ws.Range("I22").Value = 1467
$ws.Range("K22").Value = 1467
$ws.Range("M22").Value = -1172
$ws.Range("H27").Value = 3326.8
$ws.Range("I27").Value = 1467
$ws.Range("K27").Value = 1467
$ws.Range("M27").Value = -1360
$ws.Range("H40").Value = 4009.2693
$ws.Range("I40").Value = 3005.5
$ws.Range("K40").Value = 3005.5
$ws.Range("M40").Value = -2869.5
$ws.Range("H55").Value = 300.48276
$ws.Range("I55").Value = 163
$ws.Range("K55").Value = 163
$ws.Range("M55").Value = 10
$ws.Range("H61").Value = 6119.2144
$ws.Range("I61").Value = 3981.2856
$ws.Range("K61").Value = 3981.2856
$ws.Range("M61").Value = -3779.2856
$ws.Range("H68").Value = 7299.8184
$ws.Range("I68").Value = 5950
$ws.Range("J68").Value = 7599.778
$ws.Range("K68").Value = 5950
$ws.Range("L68").Value = 7599.778
$ws.Range("M68").Value = -5201
$ws.Range("N68").Value = -9097.778
$ws.Range("H71").Value = 7299.8184
$ws.Range("I71").Value = 5950
$ws.Range("J71").Value = 7599.778
$ws.Range("K71").Value = 29750
$ws.Range("L71").Value = 37998.89
$ws.Range("M71").Value = -26006
$ws.Range("N71").Value = -45486.89
$ws.Range("H95").Value = 39999
$ws.Range("J95").Value = 39999
$ws.Range("L95").Value = 39999
$ws.Range("N95").Value = -45491
$ws.Range("H100").Value = 3411.5833
$ws.Range("J100").Value = 2959.8
$ws.Range("L100").Value = 2959.8
$ws.Range("N100").Value = -4041.8
$ws.Range("H113").Value = 6119.2144
$ws.Range("I113").Value = 3981.2856
$ws.Range("K113").Value = 3981.2856
$ws.Range("M113").Value = -1811.2856
$ws.Range("H122").Value = 3857.025
$ws.Range("I122").Value = 2976.8147
$ws.Range("K122").Value = 8930.444100000001
$ws.Range("M122").Value = -6480.444100000001
$ws.Range("H126").Value = 3156.0667
$ws.Range("I126").Value = 2535.1365
$ws.Range("J126").Value = 4863.625
$ws.Range("K126").Value = 7605.4095
$ws.Range("L126").Value = 14590.875
$ws.Range("M126").Value = -5135.4095
$ws.Range("N126").Value = -19530.875
$ws.Range("H129").Value = 84974
$ws.Range("J129").Value = 84974
$ws.Range("L129").Value = 84974
$ws.Range("N129").Value = -94974
$ws.Range("H132").Value = 10210224
$ws.Range("I132").Value = 20002640
$ws.Range("K132").Value = 60007920
$ws.Range("M132").Value = -60005390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 557.9231
$ws.Range("I100").Value = 354.05884
$ws.Range("J100").Value = 943
$ws.Range("K100").Value = 708.11768
$ws.Range("L100").Value = 1886
$ws.Range("M100").Value = -167.11768
$ws.Range("N100").Value = -2968
$ws.Range("H107").Value = 13889891
$ws.Range("I107").Value = 485.26666
$ws.Range("K107").Value = 1455.79998
$ws.Range("M107").Value = 464.20002
$ws.Range("H108").Value = 69999
$ws.Range("J108").Value = 69999
$ws.Range("L108").Value = 69999
$ws.Range("N108").Value = -77679
$ws.Range("H113").Value = 12993.954
$ws.Range("J113").Value = 1861.6
$ws.Range("L113").Value = 5584.799999999999
$ws.Range("N113").Value = -9924.799999999999
$ws.Range("H122").Value = 579144.1
$ws.Range("I122").Value = 1339334.6
$ws.Range("K122").Value = 4018003.8
$ws.Range("M122").Value = -4015553.8
$ws.Range("H126").Value = 567.2353000000001
$ws.Range("I126").Value = 513
$ws.Range("J126").Value = 666.6667
$ws.Range("K126").Value = 1539
$ws.Range("L126").Value = 2000.0001
$ws.Range("M126").Value = 931
$ws.Range("N126").Value = -6940.0001
$ws.Range("H132").Value = 6806.1665
$ws.Range("I132").Value = 7398.7896
$ws.Range("K132").Value = 22196.3688
$ws.Range("M132").Value = -19666.3688
$ws.Range("H136").Value = 50507196
$ws.Range("I136").Value = 200002980
$ws.Range("K136").Value = 600008940
$ws.Range("M136").Value = -600006390
